# Update column G ("K" - strikeouts) values on Sheet1 with newly regenerated
# simulation results, per the commit "regen save_data to use K instead of
# Strike#, regen std/mean, calc and write s_vals".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newK = @{
    2  = 0
    3  = 1
    4  = 0
    5  = 1
    6  = 0
    7  = 1
    8  = 1
    9  = 1
    10 = 0
    11 = 1
    12 = 1
    13 = 2
    14 = 2
    15 = 0
    16 = 2
    17 = 2
    18 = 0
    19 = 1
    20 = 5
    21 = 3
    22 = 0
    23 = 2
    24 = 0
    25 = 2
    26 = 1
    27 = 2
    28 = 1
    29 = 1
    30 = 0
    31 = 2
    32 = 1
    33 = 2
    34 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
